$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shp = $master.Shapes.AddTextbox(1, 232.84614173228346, 13.123464566929133, 494.3077165354331, 20.599212598425197)
Write-Host $shp.Name
